# Actualizacion automatica 2025-06-23 11:00:08
#
# A new client, "WONG SANCHEZ CLAUDIA PAULINA", is inserted for advisor
# "GUERRERO FAREZ FABIAN MAURICIO" right before the existing
# "ZUÑIGA CORONEL MARCIA LUZMILA" row (row 52) on both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets. Inserting the row pushes
# the old row 52 down to row 53 and the trailing totals/summary row down
# to row 54, whose "N de 51" counters are refreshed to "N de 52" (one more
# client in the denominator). "CUMPLIMIENTO MENSUAL" is untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO  (columns C..R, dimension A1:R53 -> A1:R54)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a blank row at 52; Excel shifts rows 52-53 down to 53-54 and
# copies the formatting of the row below into the new blank row 52.
$ws1.Rows.Item(52).Insert()

# New row 52 only carries the client name (advisor column A stays blank,
# matching the repeated-advisor layout already used lower in the sheet).
$ws1.Range("B52").Value = "WONG SANCHEZ CLAUDIA PAULINA"
$ws1.Range("C52:R52").Value = 0

# Row 54 is the former summary row ("N de 51"); bump every "de 51" counter
# to "de 52" now that there is one more client row feeding it.
$cols1 = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
$nums1 = @(6,6,2,0,0,0,1,1,2,4,7,0,3,4,1,0)
for ($i = 0; $i -lt $cols1.Length; $i++) {
    $ws1.Range($cols1[$i] + "54").Value = ($nums1[$i].ToString() + " de 52")
}

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL  (columns C..G, dimension A1:G53 -> A1:G54)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(52).Insert()

$ws2.Range("B52").Value = "WONG SANCHEZ CLAUDIA PAULINA"
$ws2.Range("C52:G52").Value = 0

# Row 54's monetary totals (96431.17, 62096, 97690.82, 46695.73, 74800)
# are simply the old row 53 shifted down by the Insert() above, so they
# need no further edits here.

Write-Output "Inserted WONG SANCHEZ CLAUDIA PAULINA row on both sheets."
